# Week1Java.pptx edit: add a new 4th slide ("What are we building?") that
# recaps the pirate-themed running example (ships / role / treasure).
$p = $ppt.ActivePresentation

# New slide goes at the end of the deck, using the same "Title and Content"
# layout (index 2) already used by slide 2 and slide 3.
$s = $p.Slides.Add($p.Slides.Count + 1, 2)

# --- Title placeholder -----------------------------------------------
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Text = "What are we building? "

# --- Body / content placeholder --------------------------------------
# Build the bullet list paragraph-by-paragraph (via InsertAfter) rather
# than one multi-line Text assignment, so every paragraph gets its own
# run/rPr instead of only the first one.
$body = $s.Shapes.Item(2).TextFrame.TextRange
$body.Text = "Pirates!"
[void]$body.InsertAfter("`rShips")
[void]$body.InsertAfter("`rRole")
[void]$body.InsertAfter("`rTreasure")
[void]$body.InsertAfter("`r")
[void]$body.InsertAfter("`r")

# Indent the three sub-bullets (Ships / Role / Treasure) one level in,
# matching <a:pPr lvl="1"/> in the slide XML. IndentLevel is 1-based, so
# level 2 == lvl="1".
$body.Paragraphs(2).IndentLevel = 2
$body.Paragraphs(3).IndentLevel = 2
$body.Paragraphs(4).IndentLevel = 2
